$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Label" header in column H, matching style of existing headers (e.g. G1) ---
$ws.Cells.Item(1, 8).Value2 = "Label"
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)  # xlPasteFormats

# --- Update refit numeric results for the "100" iterations block (rows 3,5,6,7,8) ---
$ws.Cells.Item(3, 4).Value2 = 0.5790346104477224
$ws.Cells.Item(3, 5).Value2 = 0.5790346104477224

$ws.Cells.Item(5, 4).Value2 = 0.5132308771104042
$ws.Cells.Item(5, 5).Value2 = 0.4867691228895958

$ws.Cells.Item(6, 4).Value2 = 0.4337411144883669
$ws.Cells.Item(6, 5).Value2 = 0.5662588855116331

$ws.Cells.Item(7, 4).Value2 = 0.5617969412088493
$ws.Cells.Item(7, 5).Value2 = 0.4382030587911507

$ws.Cells.Item(8, 4).Value2 = 0.4493106434521933
$ws.Cells.Item(8, 5).Value2 = 0.5506893565478066
$ws.Cells.Item(8, 6).Value2 = 0.73952716588974

# --- Populate new "Label" column (H) for all data rows ---
$ws.Cells.Item(2, 8).Value2 = 0
$ws.Cells.Item(3, 8).Value2 = 0
$ws.Cells.Item(4, 8).Value2 = 1
$ws.Cells.Item(5, 8).Value2 = 1
$ws.Cells.Item(6, 8).Value2 = 1
$ws.Cells.Item(7, 8).Value2 = 1
$ws.Cells.Item(8, 8).Value2 = 1

$ws.Cells.Item(9, 8).Value2 = 0
$ws.Cells.Item(10, 8).Value2 = 0
$ws.Cells.Item(11, 8).Value2 = 1
$ws.Cells.Item(12, 8).Value2 = 1
$ws.Cells.Item(13, 8).Value2 = 1
$ws.Cells.Item(14, 8).Value2 = 1
$ws.Cells.Item(15, 8).Value2 = 1

Write-Host "Edit complete"
